{"js": "// Replace multiplication problems in the table with the new values from\n// the commit, matching old text (no duplicates in this document) and\n// inserting the new text in its place while preserving run formatting.\nconst replacements = [\n  [\"261\u00d76=\", \"160\u00d77=\"],\n  [\"101\u00d76=\", \"110\u00d74=\"],\n  [\"677\u00d74=\", \"362\u00d72=\"],\n  [\"661\u00d76=\", \"535\u00d79=\"],\n  [\"276\u00d77=\", \"372\u00d75=\"],\n  [\"220\u00d76=\", \"854\u00d72=\"],\n  [\"658\u00d72=\", \"266\u00d74=\"],\n  [\"607\u00d75=\", \"247\u00d76=\"],\n  [\"953\u00d74=\", \"135\u00d75=\"],\n  [\"967\u00d75=\", \"913\u00d79=\"],\n  [\"268\u00d75=\", \"944\u00d73=\"],\n  [\"743\u00d75=\", \"272\u00d73=\"],\n  [\"864\u00d74=\", \"355\u00d75=\"],\n  [\"642\u00d78=\", \"475\u00d72=\"],\n  [\"899\u00d76=\", \"432\u00d79=\"],\n  [\"320\u00d73=\", \"532\u00d75=\"],\n  [\"653\u00d75=\", \"429\u00d79=\"],\n  [\"137\u00d72=\", \"110\u00d78=\"],\n  [\"633\u00d77=\", \"947\u00d75=\"],\n  [\"947\u00d78=\", \"270\u00d72=\"],\n  [\"959\u00d79=\", \"708\u00d78=\"],\n  [\"767\u00d75=\", \"350\u00d79=\"],\n  [\"775\u00d72=\", \"203\u00d78=\"],\n  [\"419\u00d77=\", \"298\u00d78=\"],\n  [\"372\u00d78=\", \"218\u00d76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace multiplication problems in the table with the new values from\n# the commit. Each old value is unique in the document, so a simple\n# Find/Replace (ReplaceAll) per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$oldTexts = @(\n    \"261\u00d76=\",\n    \"101\u00d76=\",\n    \"677\u00d74=\",\n    \"661\u00d76=\",\n    \"276\u00d77=\",\n    \"220\u00d76=\",\n    \"658\u00d72=\",\n    \"607\u00d75=\",\n    \"953\u00d74=\",\n    \"967\u00d75=\",\n    \"268\u00d75=\",\n    \"743\u00d75=\",\n    \"864\u00d74=\",\n    \"642\u00d78=\",\n    \"899\u00d76=\",\n    \"320\u00d73=\",\n    \"653\u00d75=\",\n    \"137\u00d72=\",\n    \"633\u00d77=\",\n    \"947\u00d78=\",\n    \"959\u00d79=\",\n    \"767\u00d75=\",\n    \"775\u00d72=\",\n    \"419\u00d77=\",\n    \"372\u00d78=\"\n)\n\n$newTexts = @(\n    \"160\u00d77=\",\n    \"110\u00d74=\",\n    \"362\u00d72=\",\n    \"535\u00d79=\",\n    \"372\u00d75=\",\n    \"854\u00d72=\",\n    \"266\u00d74=\",\n    \"247\u00d76=\",\n    \"135\u00d75=\",\n    \"913\u00d79=\",\n    \"944\u00d73=\",\n    \"272\u00d73=\",\n    \"355\u00d75=\",\n    \"475\u00d72=\",\n    \"432\u00d79=\",\n    \"532\u00d75=\",\n    \"429\u00d79=\",\n    \"110\u00d78=\",\n    \"947\u00d75=\",\n    \"270\u00d72=\",\n    \"708\u00d78=\",\n    \"350\u00d79=\",\n    \"203\u00d78=\",\n    \"298\u00d78=\",\n    \"218\u00d76=\"\n)\n\nfor ($i = 0; $i -lt $oldTexts.Count; $i++) {\n    $old = $oldTexts[$i]\n    $new = $newTexts[$i]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.Text = $old\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
